$wb = $excel.ActiveWorkbook

# --- Swap the names of "Acc_Upfront1" and "Acc_Upfront3" tabs -------------
# (the underlying sheets/data stay put - only the tab captions swap)
$wb.Worksheets("Acc_Upfront1").Name = "Acc_UpfrontTEMP"
$wb.Worksheets("Acc_Upfront3").Name = "Acc_Upfront1"
$wb.Worksheets("Acc_UpfrontTEMP").Name = "Acc_Upfront3"

# --- Transactions sheet: move the selected cell to D5 ---------------------
$wsTrans = $wb.Worksheets("Transactions")
$wsTrans.Range("D5").Select() | Out-Null

# --- Acc_Upfront2: widen column G to fit its contents ----------------------
$wsUp2 = $wb.Worksheets("Acc_Upfront2")
$wsUp2.Columns.Item(7).AutoFit() | Out-Null

# --- Sheet now named "Acc_Upfront1" (originally "Acc_Upfront3"): ----------
# move the selected cell to I20; this tab ends up NOT active
$wsUp1 = $wb.Worksheets("Acc_Upfront1")
$wsUp1.Range("I20").Select() | Out-Null

# --- Sheet now named "Acc_Upfront3" (originally "Acc_Upfront1"): ----------
# widen column G, select H23, and leave this as the final active tab
$wsUp3 = $wb.Worksheets("Acc_Upfront3")
$wsUp3.Columns.Item(7).AutoFit() | Out-Null
$wsUp3.Activate() | Out-Null
$wsUp3.Range("H23").Select() | Out-Null
